# "Generate Report for Handback" — populate the Latest Target File /
# Latest Handback File / Latest Handback DateTime columns for both the
# zh-cn and de-de handback reports, and flip the status text from
# "Ready for handoff" to "Handed back: in sync with en-US" everywhere
# it appears (Overview summary + per-language Status column).

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

$mdFile1 = "054ce818-756e-4097-9fbd-6bb3db15773f.md"
$mdFile2 = "4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bfdbf6a5165793a6dae619fb53c9a037e694291/e2e/054ce818-756e-4097-9fbd-6bb3db15773f.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bfdbf6a5165793a6dae619fb53c9a037e694291/e2e/4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.md"

# ----- Overview sheet: status text for both languages, both rows -----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

# ----- zh-cn sheet -----
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("C3").Value = $handedBack

# widen the Status / Latest Target File / Latest Handback File columns
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# row 2 - 054ce818 file
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl1, "", "", $mdFile1)
$zhcn.Range("J2").Value = "054ce818-756e-4097-9fbd-6bb3db15773f.c43f06c6884a0092f01462099fae516026448794.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-15 14:45:43"

# row 3 - 4f38d11f file
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdUrl2, "", "", $mdFile2)
$zhcn.Range("J3").Value = "4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.5cc6957cee5a60b991cd9ec994a5419a2287c524.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-15 14:45:43"

# ----- de-de sheet -----
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $handedBack
$dede.Range("C3").Value = $handedBack

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

# row 2 - 054ce818 file
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl1, "", "", $mdFile1)
$dede.Range("J2").Value = "054ce818-756e-4097-9fbd-6bb3db15773f.c43f06c6884a0092f01462099fae516026448794.de-de.xlf"
$dede.Range("K2").Value = "2016-08-15 14:45:51"

# row 3 - 4f38d11f file
$dede.Hyperlinks.Add($dede.Range("I3"), $mdUrl2, "", "", $mdFile2)
$dede.Range("J3").Value = "4f38d11f-f14e-400f-b6ba-7c0ecfb544fc.5cc6957cee5a60b991cd9ec994a5419a2287c524.de-de.xlf"
$dede.Range("K3").Value = "2016-08-15 14:45:51"
